# Delivery address sheet added ("delivery address codea added")
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Keep Sheet1's original selection info consistent with the final file
# (Excel records the last selection on a sheet when focus leaves it).
$ws1.Range("A1:B1").Select()

# Insert the new "Sheet2" right after Sheet1, then make it the active tab.
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"

# ---- Header row (row 1): bold, yellow fill, thin box border ----
$headerRange = $ws2.Range("A1:H1")
$ws1.Range("A1").Copy()
$headerRange.PasteSpecial(-4122)
$headerRange.Borders.LineStyle = 1

$ws2.Range("A1").Value = "Full Name"
$ws2.Range("B1").Value = "Address 1"
$ws2.Range("C1").Value = "Address 2"
$ws2.Range("D1").Value = "Landmark"
$ws2.Range("E1").Value = "City"
$ws2.Range("F1").Value = "State"
$ws2.Range("G1").Value = "Pin Code"
$ws2.Range("H1").Value = "Mobile Number"

# ---- Data row (row 2): thin box border, default style ----
$dataRange = $ws2.Range("A2:H2")
$dataRange.Borders.LineStyle = 1

$ws2.Range("A2").Value = "Valmiki"
$ws2.Range("B2").Value = "203, Akshar Bluechip IT Park,"
$ws2.Range("C2").Value = "Turbhe MIDC"
$ws2.Range("D2").Value = "Turbhe"
$ws2.Range("E2").Value = "New Mumbai"
$ws2.Range("F2").Value = "Maharashtra "
$ws2.Range("G2").Value = 400705
$ws2.Range("H2").Value = 2262596124

# ---- Column widths (best-fit to content) ----
$ws2.Columns.Item(1).ColumnWidth = 9.230769230769232
$ws2.Columns.Item(2).ColumnWidth = 25.59031249999998
$ws2.Columns.Item(3).ColumnWidth = 11.585781249999975
$ws2.Columns.Item(4).ColumnWidth = 8.590312499999978
$ws2.Columns.Item(5).ColumnWidth = 11.935468749999977
$ws2.Columns.Item(6).ColumnWidth = 11.585781249999975
$ws2.Columns.Item(7).ColumnWidth = 7.935468749999977
$ws2.Columns.Item(8).ColumnWidth = 14.425156249999976

# Land on I2 and make Sheet2 the active/visible tab, matching the saved view.
$ws2.Range("I2").Select()
$ws2.Activate()
